$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 56.426
$ws.Range("D2").Value = 56.426
$ws.Range("E2").Value = 3.03340581
$ws.Range("F2").Value = 0.04028454
$ws.Range("G2").Value = 2.26997876
$ws.Range("H2").Value = 129.56059548
$ws.Range("I2").Value = 6.449971017464042
$ws.Range("J2").Value = 6.449971017464042
$ws.Range("K2").Value = 0.3402677275178092
$ws.Range("L2").Value = 0.003814876311188826
$ws.Range("M2").Value = 0.318706954806613
$ws.Range("N2").Value = 31.02704578393793
$ws.Range("C3").Value = 86.163
$ws.Range("D3").Value = 86.163
$ws.Range("E3").Value = 1.99537222
$ws.Range("F3").Value = 0.03343293999999999
$ws.Range("G3").Value = 2.84798751
$ws.Range("H3").Value = 246.73747622
$ws.Range("I3").Value = 11.32257177639566
$ws.Range("J3").Value = 11.32257177639566
$ws.Range("K3").Value = 0.2629915256928592
$ws.Range("L3").Value = 0.004759148206144462
$ws.Range("M3").Value = 0.3532168311495181
$ws.Range("N3").Value = 51.72210462156272
$ws.Range("C4").Value = 31.54
$ws.Range("D4").Value = 63.035
$ws.Range("E4").Value = 2.729127
$ws.Range("F4").Value = 0.03505555
$ws.Range("G4").Value = 0.5526313900000001
$ws.Range("H4").Value = 17.74162784
$ws.Range("I4").Value = 4.284244945699855
$ws.Range("J4").Value = 8.544634912375471
$ws.Range("K4").Value = 0.3610348958957492
$ws.Range("L4").Value = 0.005039072190993699
$ws.Range("M4").Value = 0.1103469623427223
$ws.Range("N4").Value = 5.601350248477205
$ws.Range("C5").Value = 47.451
$ws.Range("D5").Value = 92.688
$ws.Range("E5").Value = 1.86378554
$ws.Range("F5").Value = 0.02807004
$ws.Range("G5").Value = 0.65785667
$ws.Range("H5").Value = 31.62002765
$ws.Range("I5").Value = 7.578194245862798
$ws.Range("J5").Value = 13.81243785287157
$ws.Range("K5").Value = 0.2769510662525184
$ws.Range("L5").Value = 0.004653765500294544
$ws.Range("M5").Value = 0.1130795881397919
$ws.Range("N5").Value = 9.24464776399464
$ws.Range("C6").Value = 18.105
$ws.Range("D6").Value = 72.288
$ws.Range("E6").Value = 2.42129739
$ws.Range("F6").Value = 0.02495465
$ws.Range("G6").Value = 0.11301532
$ws.Range("H6").Value = 2.11541271
$ws.Range("I6").Value = 3.242914043185772
$ws.Range("J6").Value = 12.93384476249581
$ws.Range("K6").Value = 0.4837121954432289
$ws.Range("L6").Value = 0.004824932357012868
$ws.Range("M6").Value = 0.03072817692207845
$ws.Range("N6").Value = 0.8867501776835333
$ws.Range("C7").Value = 26.451
$ws.Range("D7").Value = 96.639
$ws.Range("E7").Value = 1.78548637
$ws.Range("F7").Value = 0.02273616
$ws.Range("G7").Value = 0.1491138
$ws.Range("H7").Value = 4.053499850000001
$ws.Range("I7").Value = 5.020637987926017
$ws.Range("J7").Value = 14.16408487337958
$ws.Range("K7").Value = 0.256407048601429
$ws.Range("L7").Value = 0.004544597564297792
$ws.Range("M7").Value = 0.03660547286682422
$ws.Range("N7").Value = 1.774858587307004
$ws.Range("C8").Value = 12.243
$ws.Range("D8").Value = 73.29000000000001
$ws.Range("E8").Value = 2.42486441
$ws.Range("F8").Value = 0.0206359
$ws.Range("G8").Value = 0.04246956
$ws.Range("H8").Value = 0.54903085
$ws.Range("I8").Value = 2.657633160166018
$ws.Range("J8").Value = 15.91720004451968
$ws.Range("K8").Value = 0.5720805643859963
$ws.Range("L8").Value = 0.004600885386133449
$ws.Range("M8").Value = 0.01475908414405078
$ws.Range("N8").Value = 0.2935980370229359
$ws.Range("C9").Value = 18.781
$ws.Range("D9").Value = 94.43899999999999
$ws.Range("E9").Value = 1.84088512
$ws.Range("F9").Value = 0.01839505
$ws.Range("G9").Value = 0.05737918
$ws.Range("H9").Value = 1.12824401
$ws.Range("I9").Value = 4.246772194746744
$ws.Range("J9").Value = 15.75063050678802
$ws.Range("K9").Value = 0.3156188089121507
$ws.Range("L9").Value = 0.003265759633451838
$ws.Range("M9").Value = 0.01606636979426174
$ws.Range("N9").Value = 0.5714073840994911
$ws.Range("C10").Value = 9.012
$ws.Range("D10").Value = 71.884
$ws.Range("E10").Value = 2.507227299999999
$ws.Range("F10").Value = 0.01660387
$ws.Range("G10").Value = 0.01896527
$ws.Range("H10").Value = 0.18379979
$ws.Range("I10").Value = 2.235811900597824
$ws.Range("J10").Value = 17.84323190806321
$ws.Range("K10").Value = 0.661523204476586
$ws.Range("L10").Value = 0.003252455070725357
$ws.Range("M10").Value = 0.006913931122789022
$ws.Range("N10").Value = 0.1091010990103269
$ws.Range("C11").Value = 14.56
$ws.Range("D11").Value = 88.339
$ws.Range("E11").Value = 1.97984356
$ws.Range("F11").Value = 0.01576841
$ws.Range("G11").Value = 0.02877109
$ws.Range("H11").Value = 0.44888761
$ws.Range("I11").Value = 3.83863714554204
$ws.Range("J11").Value = 16.40900162515484
$ws.Range("K11").Value = 0.3711075514913437
$ws.Range("L11").Value = 0.002775071313082367
$ws.Range("M11").Value = 0.009467816879744156
$ws.Range("N11").Value = 0.2841402757676236
$ws.Range("C12").Value = 7.083
$ws.Range("D12").Value = 70.577
$ws.Range("E12").Value = 2.58359017
$ws.Range("F12").Value = 0.01505144
$ws.Range("G12").Value = 0.01093338
$ws.Range("H12").Value = 0.08521947000000001
$ws.Range("I12").Value = 1.940070842746438
$ws.Range("J12").Value = 19.28839322268999
$ws.Range("K12").Value = 0.7364057312791176
$ws.Range("L12").Value = 0.003110658857046996
$ws.Range("M12").Value = 0.004576088937680907
$ws.Range("N12").Value = 0.05857671725885571
$ws.Range("C13").Value = 11.58
$ws.Range("D13").Value = 80.089
$ws.Range("E13").Value = 2.19018567
$ws.Range("F13").Value = 0.01445044
$ws.Range("G13").Value = 0.01685115000000001
$ws.Range("H13").Value = 0.21177302
$ws.Range("I13").Value = 3.202163933007467
$ws.Range("J13").Value = 15.70362005025219
$ws.Range("K13").Value = 0.4227209638447388
$ws.Range("L13").Value = 0.00265218383327504
$ws.Range("M13").Value = 0.006173582420070003
$ws.Range("N13").Value = 0.1475696315894535
